# Generate Report for Handoff
# Updates the localization status report: marks zh-cn/de-de as
# "Ready for handoff" (was "In Translation") and refreshes the
# handoff timestamps, then widens the status-datetime columns so the
# new values are fully visible.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Handoff timestamps -----------------------------------------------
# zh-cn handoff time moves from 09:02:08 -> 09:02:44
$wsZhCn.Range("H2").Value = "2016-09-04 09:02:44"

# de-de handoff time (and the workbook-wide "latest" generate date on the
# Overview sheet) moves from 09:02:13 -> 09:02:49
$wsDeDe.Range("H2").Value     = "2016-09-04 09:02:49"
$wsOverview.Range("G2").Value = "2016-09-04 09:02:49"

# --- Column widths ------------------------------------------------------
# Widen the status/datetime columns to fit the new, longer text.
$newWidth = 17.2159881591797 - (5.0 / 6.0)

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth   # column F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth     = $newWidth   # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = $newWidth   # column C (Status)
